# Weekly cryptos data refresh (prices in column D, 1h volume change in column E).
# A handful of rows also had their coin re-ranked (Fetch.AI/NEARProtocol/FirstDigitalUSD
# around rows 37-39, and Stellar/VeChain swapped around rows 48-49), so B/C get rewritten too.
#
# Column D often holds values that *look* numeric ("0.999", "142.00", ...) but must stay
# as literal text (matches the sheet's existing inlineStr cells, and preserves formatting
# like trailing zeros that a real number would drop). A leading single-quote forces Excel
# to store the value as text (quoted-text / "Text" cell) instead of coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.199.25'
$ws.Range("E2").Value = '  -15.51%  '

$ws.Range("D3").Value = '2.247.19'
$ws.Range("E3").Value = '  -22.55%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''432.86'
$ws.Range("E5").Value = '  -17.77%  '

$ws.Range("D6").Value = '''116.83'
$ws.Range("E6").Value = '  -18.30%  '

$ws.Range("D7").Value = '''0.996'
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").Value = '''0.454'
$ws.Range("E8").Value = '  -16.92%  '

$ws.Range("D9").Value = '2.226.72'
$ws.Range("E9").Value = '  -23.45%  '

$ws.Range("D10").Value = '''5.11'
$ws.Range("E10").Value = '  -14.49%  '

$ws.Range("D11").Value = '''0.0841'
$ws.Range("E11").Value = '  -21.43%  '

$ws.Range("D12").Value = '''0.295'
$ws.Range("E12").Value = '  -17.61%  '

$ws.Range("D13").Value = '''0.120'
$ws.Range("E13").Value = '  -7.07%  '

$ws.Range("D14").Value = '2.613.67'
$ws.Range("E14").Value = '  -23.23%  '

$ws.Range("D15").Value = '51.105.92'
$ws.Range("E15").Value = '  -15.64%  '

$ws.Range("D16").Value = '''18.37'
$ws.Range("E16").Value = '  -18.50%  '

$ws.Range("E17").Value = '  -19.35%  '

$ws.Range("D18").Value = '2.243.62'
$ws.Range("E18").Value = '  -22.81%  '

$ws.Range("D19").Value = '''3.86'
$ws.Range("E19").Value = '  -22.07%  '

$ws.Range("D20").Value = '''292.25'
$ws.Range("E20").Value = '  -16.55%  '

$ws.Range("D21").Value = '''0.997'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").Value = '''5.68'
$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("D23").Value = '''8.51'
$ws.Range("E23").Value = '  -26.36%  '

$ws.Range("D24").Value = '''4.95'
$ws.Range("E24").Value = '  -23.87%  '

$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("D26").Value = '''52.40'
$ws.Range("E26").Value = '  -19.09%  '

$ws.Range("D27").Value = '''0.359'
$ws.Range("E27").Value = '  -20.25%  '

$ws.Range("D28").Value = '2.347.38'
$ws.Range("E28").Value = '  -22.49%  '

$ws.Range("D29").Value = '''0.135'
$ws.Range("E29").Value = '  -23.93%  '

$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.10%  '

$ws.Range("D31").Value = '''6.68'
$ws.Range("E31").Value = '  -14.51%  '

$ws.Range("D32").Value = '''142.00'
$ws.Range("E32").Value = '  -6.51%  '

$ws.Range("D33").Value = '0.0₃0622'
$ws.Range("E33").Value = '  -27.08%  '

$ws.Range("D34").Value = '''16.29'
$ws.Range("E34").Value = '  -16.51%  '

$ws.Range("D35").Value = '''1.29'
$ws.Range("E35").Value = '  -22.62%  '

$ws.Range("D36").Value = '''4.56'
$ws.Range("E36").Value = '  -17.78%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '''0.774'
$ws.Range("E37").Value = '  -22.03%  '

$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''3.28'
$ws.Range("E38").Value = '  -23.32%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''0.989'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").Value = '''0.965'
$ws.Range("E40").Value = '  -19.03%  '

$ws.Range("D41").Value = '''31.43'
$ws.Range("E41").Value = '  -16.37%  '

$ws.Range("D42").Value = '''10.09'
$ws.Range("E42").Value = '  -2.29%  '

$ws.Range("D43").Value = '''0.540'
$ws.Range("E43").Value = '  -16.70%  '

$ws.Range("D44").Value = '''0.0489'
$ws.Range("E44").Value = '  -15.49%  '

$ws.Range("D45").Value = '''3.02'
$ws.Range("E45").Value = '  -18.47%  '

$ws.Range("D46").Value = '1.843.99'
$ws.Range("E46").Value = '  -19.42%  '

$ws.Range("D47").Value = '''1.13'
$ws.Range("E47").Value = '  -22.56%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0199'
$ws.Range("E48").Value = '  -15.73%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.0798'
$ws.Range("E49").Value = '  -12.88%  '

$ws.Range("D50").Value = '''3.93'
$ws.Range("E50").Value = '  -20.29%  '

$ws.Range("D51").Value = '''4.62'
$ws.Range("E51").Value = '  -5.15%  '
